# Update "想去人数" (interest count) values in column F across sheets,
# reflecting newly generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 640
$ws1.Range("F5").Value = 561
$ws1.Range("F6").Value = 309
$ws1.Range("F7").Value = 2778
$ws1.Range("F9").Value = 7803
$ws1.Range("F10").Value = 203
$ws1.Range("F13").Value = 342
$ws1.Range("F14").Value = 48

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 5

# --- Sheet "全部类型" (all types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 640
$ws4.Range("F5").Value = 561
$ws4.Range("F6").Value = 309
$ws4.Range("F9").Value = 2778
$ws4.Range("F11").Value = 7803
$ws4.Range("F12").Value = 203
$ws4.Range("F15").Value = 5
$ws4.Range("F17").Value = 342
$ws4.Range("F18").Value = 48
